$wb = $excel.ActiveWorkbook

# Current layout: Sheet1, Sheet2 (drying data + chart).
# Target layout:  Sheet1, Sheet3 (new weigh/dry-weigh data), Sheet2 (unchanged drying data).
# Insert the new worksheet right after Sheet1 - Excel names it "Sheet3"
# automatically (the next unused default sheet name) and makes it active.
$sheet1 = $wb.Worksheets.Item("Sheet1")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)

# Populate the new "Sheet3" tab with the measured drying data.
$newSheet.Range("D2").Value = "Dry weigh"
$newSheet.Range("D3").Value = 36.7

$gVals = @(2.46, 2.47, 2.56, 2.73, 2.92, 3.05, 3.32, 3.54, 3.7, 3.91, 4.2)
$hVals = @(100.9, 98.6, 96.6, 89.5, 84.4, 80.3, 74.6, 69.7, 65.2, 59.8, 54.7)

for ($i = 0; $i -lt $gVals.Length; $i++) {
    $row = 3 + $i
    $newSheet.Cells.Item($row, 7).Value = $gVals[$i]
    $newSheet.Cells.Item($row, 8).Value = $hVals[$i]
}

# Match the author's final selection on the new sheet.
$newSheet.Range("H13").Select() | Out-Null
